$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The new shape is a small right-aligned "t" type-annotation label, identical
# in style to the other such labels already on the diagram (e.g. "TextBox 92"
# sitting next to the first "result" callout). Duplicate that shape so the
# new one inherits the exact same text formatting (bodyPr/lstStyle: 14pt,
# 0070C0, right aligned), then move/rename/retext the copy into place next
# to the TagCommandParser return arrow.
$styleSource = $s.Shapes.Item("TextBox 92")
$newShapes = $styleSource.Duplicate()
$newShape = $newShapes.Item(1)

$newShape.Name = "TextBox 176"
$newShape.Left = 375.283
$newShape.Top = 569.442
$newShape.Width = 17.3499
$newShape.Height = 16.9641

$newShape.TextFrame.TextRange.Text = "t"
